$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price/Volume cells hold text like "1.015" or "28.220.00" that Excel would
# otherwise auto-convert to numbers; force text via NumberFormat, then restore
# the default "Normal" style so no stray per-cell formatting is introduced.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.220.00'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.99%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.931.82'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.66%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.015'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.96%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '320.70'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.21%  '
$ws.Range("E6").Value = '  +0.74%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4731'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -5.28%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4045'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.99%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '53.33'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.14%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08466'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -7.74%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.046'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.81%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.18'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.63%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.971.22'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.82%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.505'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.49%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.100'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -5.31%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.016'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.03%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '89.59'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.12%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001067'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.78%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06623'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.28%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.20'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.56%  '
$ws.Range("E21").Value = '  +0.71%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.795'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.77%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.290.80'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.85%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.38'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.75%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.296'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.50%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.205.18'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.68%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '155.16'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.78%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.12'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.37%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.159'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.62%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.737'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -7.93%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '123.55'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.31%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9757'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -6.82%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09583'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.79%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.672'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.48%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.439'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.93%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.583'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.55%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '9.135'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.08%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02312'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.72%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06151'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.31%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.239'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.84%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6183'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.08%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.11'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.78%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.012'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.71%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1905'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.21%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.324'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.42%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5897'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.64%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.88'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.47%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.040'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -7.16%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.397'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.21%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06770'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.13%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '109.71'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.58%  '
